# FEED/Reporting.xlsx - "added icons to header, added £ to merchant transactions"
#
# Functional changes applied here (matching the canonical-XML diff):
#   1. Fix typo in the "Styling" iteration comment cell (E49):
#        "understanding/Styling Home page" -> "Understanding/Styling Home page"
#   2. Append a new timesheet entry as row 51 (same iteration/feature/comment
#      text as row 50, 2 manhours, no G (cost) formula yet), which also bumps
#      the Total Hours (K8) and Total Cost (K9) formulas/results further down
#      the sheet via recalculation.
#   3. Select row 50 (mirrors the author having just finished editing it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Capitalise "Understanding" in the comment for row 49.
$ws.Cells.Item(49, 5).Value = "Understanding/Styling Home page"

# 2. Add the new row (51) duplicating the Styling / Styling Tables-Icons entry.
$ws.Cells.Item(51, 1).Value = 43508
$ws.Cells.Item(51, 1).NumberFormat = $ws.Cells.Item(50, 1).NumberFormat
$ws.Cells.Item(51, 2).Value = 2
$ws.Cells.Item(51, 3).Value = "Styling"
$ws.Cells.Item(51, 4).Value = 2
$ws.Cells.Item(51, 5).Value = "Styling Tables/Icons"
$ws.Cells.Item(51, 6).Value = "Getting a bit more confident – afraid to do the layout and break it all"

# 3. Update the active selection to row 50 (whole row), matching the saved view state.
[void]$ws.Range("A50").EntireRow.Select()
